# PL_SP17.xlsx — add processed P-vs-I fit results (alpha, pmax) to each
# species sheet: new header cells E1/F1, and constant per-sheet alpha/pmax
# values copied down rows 2:25 (E gets a scientific-notation number format
# on most sheets; Crypto keeps General format).

$wb = $excel.ActiveWorkbook

# sheet index (1-based, tab order) -> alpha, pmax, and whether the alpha
# column (E) should get the 0.00E+00 number format applied.
$sheetData = @(
    @{ Index = 1; Alpha = [double]"5.3581390000000001E-5"; Pmax = [double]"0.29945769999999999"; FormatAlpha = $true;  SelRange = "E2:F25"; ActiveCell = "E2" },
    @{ Index = 2; Alpha = [double]"6.6599279999999996E-5"; Pmax = [double]"0.24368890000000001"; FormatAlpha = $true;  SelRange = "E2:F25"; ActiveCell = "E2" },
    @{ Index = 3; Alpha = [double]"8.2464400000000006E-6"; Pmax = [double]"0.31484630000000002"; FormatAlpha = $true;  SelRange = "H21";    ActiveCell = "H21" },
    @{ Index = 4; Alpha = [double]"9.1315790000000003E-5"; Pmax = [double]"0.20604259999999999"; FormatAlpha = $true;  SelRange = "G19";    ActiveCell = "G19" },
    @{ Index = 5; Alpha = [double]"1.9536369999999999E-4"; Pmax = [double]"-5.785831E-3";         FormatAlpha = $false; SelRange = "I23";    ActiveCell = "I23" },
    @{ Index = 6; Alpha = [double]"5.480004E-5";            Pmax = [double]"0.2671887";            FormatAlpha = $true;  SelRange = "I9";     ActiveCell = "I9" }
)

foreach ($sd in $sheetData) {
    $ws = $wb.Worksheets.Item($sd.Index)

    # Headers
    $ws.Range("E1").Value = "alpha"
    $ws.Range("F1").Value = "pmax"

    # Constant alpha/pmax values copied down for every data row (2-25)
    $ws.Range("E2:E25").Value = $sd.Alpha
    $ws.Range("F2:F25").Value = $sd.Pmax

    if ($sd.FormatAlpha) {
        $ws.Range("E2:E25").NumberFormat = "0.00E+00"
    }

    # Restore the sheet's selection / active cell as saved in the source file
    $ws.Activate() | Out-Null
    $ws.Range($sd.SelRange).Select() | Out-Null
}

# Leave the workbook on the originally active sheet (Plocamium, tab 6)
$lastWs = $wb.Worksheets.Item(6)
$lastWs.Activate() | Out-Null
$lastWs.Range("I9").Select() | Out-Null
